# ToDo.xlsx edit: "Added mouse device to input manager. Implemented mouse
# controls in model viewer."
#
# The old row 2 ("Engine | Add input manager with keyboard and mouse | 6")
# is removed entirely (that task is considered done / folded into the
# mouse-controls work below), which shifts every following row up by one.
# The (new) row 2 - "Model Viewer | Add camera controls to model viewer" -
# also gets its estimate revised from 3 down to 2, and loses the leftover
# bold-ish style that the old row 2 used to carry.
#
# Cell comments are anchored to absolute cell refs in this file format, so
# after the row shift the three existing comments (originally on B5, B13,
# B18) have to be moved down to B4, B12, B17 to stay on the same logical
# tasks, and a brand new comment is added on B2 ("Add camera controls to
# model viewer").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the completed "Add input manager with keyboard and mouse" row.
# Everything below (rows 3-22) shifts up to rows 2-21.
$ws.Rows.Item(2).Delete()

# --- The row that is now row 2 ("Add camera controls to model viewer")
# gets a revised estimate and drops its old (no-op bold) cell style.
$ws.Range("C2").Value = 2
$ws.Range("A2:C2").ClearFormats()

# --- Relocate the existing comments to follow their tasks up one row.
# (Comment anchors don't auto-shift with row deletes, so move them by hand:
# grab the original text, delete the old comment, re-add at the new cell.)
$note1 = $ws.Range("B5").Comment.Text()
$ws.Range("B5").Comment.Delete()
$ws.Range("B4").AddComment($note1)

$note2 = $ws.Range("B13").Comment.Text()
$ws.Range("B13").Comment.Delete()
$ws.Range("B12").AddComment($note2)

$note3 = $ws.Range("B18").Comment.Text()
$ws.Range("B18").Comment.Delete()
$ws.Range("B17").AddComment($note3)

# --- New comment on the "Add camera controls to model viewer" row noting
# the mouse-look bug that prompted this change.
$ws.Range("B2").AddComment("Jonny:`nThe camera inadvertantly rolls around :-(")

# --- Match the author's final selection.
$ws.Range("B2").Select()
